# Update "paises" (countries) workbook:
#  - Refresh the "last updated" timestamp in A1
#  - Swap four pairs of country rows whose ranking changed (new data makes
#    one overtake the other), keeping the table sorted by "Casos totales" desc
#  - Update the numeric COVID-19 stats for the rows that changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 17:39"

# --- Row label swaps (rank changed between the two countries) ----------
# Honduras / Azerbaiyan (rows 60-61)
$ws.Cells.Item(60, 1).Value = "Azerbaiyan"
$ws.Cells.Item(61, 1).Value = "Honduras"

# Niger / Jordania (rows 128-129)
$ws.Cells.Item(128, 1).Value = "Jordania"
$ws.Cells.Item(129, 1).Value = "Niger"

# Dominica / Fiyi (rows 202-203)
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"

# Groenlandia / Islas Malvinas (rows 208-209)
$ws.Cells.Item(208, 1).Value = "Islas Malvinas"
$ws.Cells.Item(209, 1).Value = "Groenlandia"

# --- Numeric updates -----------------------------------------------------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 2434773
$ws.Cells.Item(4, 3).Value = 10605
$ws.Cells.Item(4, 4).Value = 1020499
$ws.Cells.Item(4, 5).Value = 1290606
$ws.Cells.Item(4, 7).Value = 195
$ws.Cells.Item(4, 8).Value = 123668

# Row 7 - India
$ws.Cells.Item(7, 2).Value = 465553
$ws.Cells.Item(7, 3).Value = 9438
$ws.Cells.Item(7, 4).Value = 264542
$ws.Cells.Item(7, 5).Value = 186377
$ws.Cells.Item(7, 7).Value = 151
$ws.Cells.Item(7, 8).Value = 14634

# Row 8 - Reino Unido
$ws.Cells.Item(8, 2).Value = 306862
$ws.Cells.Item(8, 3).Value = 652
$ws.Cells.Item(8, 7).Value = 154
$ws.Cells.Item(8, 8).Value = 43081

# Row 39 - Irak
$ws.Cells.Item(39, 2).Value = 36702
$ws.Cells.Item(39, 3).Value = 2200
$ws.Cells.Item(39, 4).Value = 16814
$ws.Cells.Item(39, 5).Value = 18558
$ws.Cells.Item(39, 7).Value = 79
$ws.Cells.Item(39, 8).Value = 1330

# Row 45 - Republica Dominicana
$ws.Cells.Item(45, 2).Value = 28631
$ws.Cells.Item(45, 3).Value = 695
$ws.Cells.Item(45, 4).Value = 16006
$ws.Cells.Item(45, 5).Value = 11934
$ws.Cells.Item(45, 7).Value = 16
$ws.Cells.Item(45, 8).Value = 691

# Row 57 - Moldavia
$ws.Cells.Item(57, 4).Value = 8400
$ws.Cells.Item(57, 5).Value = 5822
$ws.Cells.Item(57, 7).Value = 2
$ws.Cells.Item(57, 8).Value = 492

# Row 60 - now Azerbaiyan
$ws.Cells.Item(60, 2).Value = 14305
$ws.Cells.Item(60, 3).Value = 590
$ws.Cells.Item(60, 4).Value = 7768
$ws.Cells.Item(60, 5).Value = 6363
$ws.Cells.Item(60, 7).Value = 7
$ws.Cells.Item(60, 8).Value = 174

# Row 61 - now Honduras
$ws.Cells.Item(61, 2).Value = 13943
$ws.Cells.Item(61, 3).Value = 587
$ws.Cells.Item(61, 4).Value = 1461
$ws.Cells.Item(61, 5).Value = 12077
$ws.Cells.Item(61, 7).Value = 10
$ws.Cells.Item(61, 8).Value = 405

# Row 79 - Tayikistan
$ws.Cells.Item(79, 2).Value = 5630
$ws.Cells.Item(79, 3).Value = 63
$ws.Cells.Item(79, 4).Value = 4194
$ws.Cells.Item(79, 5).Value = 1384

# Row 82 - Kenia
$ws.Cells.Item(82, 4).Value = 1823
$ws.Cells.Item(82, 5).Value = 3253
$ws.Cells.Item(82, 7).Value = 2
$ws.Cells.Item(82, 8).Value = 130

# Row 84 - Guinea
$ws.Cells.Item(84, 2).Value = 5104
$ws.Cells.Item(84, 3).Value = 64
$ws.Cells.Item(84, 4).Value = 3744
$ws.Cells.Item(84, 5).Value = 1331
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 29

# Row 94 - Grecia
$ws.Cells.Item(94, 2).Value = 3310
$ws.Cells.Item(94, 3).Value = 8
$ws.Cells.Item(94, 5).Value = 1746

# Row 120 - Zambia
$ws.Cells.Item(120, 2).Value = 1489
$ws.Cells.Item(120, 3).Value = 12
$ws.Cells.Item(120, 4).Value = 1223
$ws.Cells.Item(120, 5).Value = 248

# Row 128 - now Jordania
$ws.Cells.Item(128, 2).Value = 1071
$ws.Cells.Item(128, 3).Value = 24
$ws.Cells.Item(128, 4).Value = 782
$ws.Cells.Item(128, 5).Value = 280
$ws.Cells.Item(128, 8).Value = 9

# Row 129 - now Niger
$ws.Cells.Item(129, 2).Value = 1051
$ws.Cells.Item(129, 4).Value = 913
$ws.Cells.Item(129, 5).Value = 71
$ws.Cells.Item(129, 8).Value = 67

# Row 158 - Mauricio
$ws.Cells.Item(158, 2).Value = 341
$ws.Cells.Item(158, 3).Value = 1
$ws.Cells.Item(158, 5).Value = 5

# Row 183 - Liechtenstein
$ws.Cells.Item(183, 2).Value = 82
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 55
$ws.Cells.Item(183, 5).Value = 26
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 1
